$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14
$ws.Range("B3").Value = 7
$ws.Range("B4").Value = 8
$ws.Range("B5").Value = 4
$ws.Range("B6").Value = 6
